# Swap the data between specific row pairs in the Artfynd sheet.
# The underlying observations for these rows were re-matched to the
# correct species/coordinates, which (since rows are otherwise mostly
# identical) shows up as a handful of swapped cell values per pair.
#
# Row pairs whose data must be exchanged: (5,6) (7,8) (13,15) (18,19)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell {
    param($row1, $row2, $col)

    $c1 = $ws.Cells.Item($row1, $col)
    $c2 = $ws.Cells.Item($row2, $col)

    $v1 = $c1.Value()
    $v2 = $c2.Value()

    if ($null -eq $v1) {
        $c2.ClearContents()
    } else {
        $c2.Value = $v1
    }

    if ($null -eq $v2) {
        $c1.ClearContents()
    } else {
        $c1.Value = $v2
    }
}

# Column numbers used below:
# A=1 B=2 D=4 E=5 F=6 G=7 H=8 Q=17 R=18 AC=29

$pairs = @(
    @(5, 6),
    @(7, 8),
    @(13, 15),
    @(18, 19)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    Swap-Cell $r1 $r2 1   # A  - Id
    Swap-Cell $r1 $r2 2   # B  - Taxonsorteringsordning
    Swap-Cell $r1 $r2 4   # D  - Rödlistade
    Swap-Cell $r1 $r2 5   # E  - TaxonId
    Swap-Cell $r1 $r2 6   # F  - Artnamn
    Swap-Cell $r1 $r2 7   # G  - Vetenskapligt namn
    Swap-Cell $r1 $r2 8   # H  - Auktor
    Swap-Cell $r1 $r2 17  # Q  - Ost
    Swap-Cell $r1 $r2 18  # R  - Nord
    Swap-Cell $r1 $r2 29  # AC - Publik kommentar
}

Write-Output "Row pairs swapped"
